$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 0.0325
$ws.Range("E2").Value = -0.04935
$ws.Range("K2").Value = 808.5999999999999
$ws.Range("L2").Value = 0.1495717800262666
$ws.Range("M2").Value = 244.9845
$ws.Range("N2").Value = 0.01145306517440148
$ws.Range("O2").Value = 0.3029736581746229
$ws.Range("P2").Value = 244.9845
$ws.Range("Q2").Value = 0.01145306517440148
$ws.Range("R2").Value = 0.3029736581746229
$ws.Range("U2").Value = 6202.900000000001
$ws.Range("V2").Value = 0.2899865827033749
$ws.Range("W2").Value = 0.06245198965912041
$ws.Range("X2").Value = 0.0726573574395516
$ws.Range("Y2").Value = -0.01020536778043118
$ws.Range("Z2").Value = 0.1739412175788706
$ws.Range("AB2").Value = 0.0498636495449136
$ws.Range("AC2").Value = -0.0498636495449136
$ws.Range("AD2").Value = 23663.4
$ws.Range("AF2").Value = 23663.4
$ws.Range("AG2").Value = 17460.5
$ws.Range("AH2").Value = 0.5252265629681914
$ws.Range("AI2").Value = 0.6064381873025066
$ws.Range("AJ2").Value = 0.4494244648758841
$ws.Range("AK2").Value = 0.5320500709989213

# --- Row 3 ---
$ws.Range("B3").Value = 'Crédit du Maroc S.A. (CBSE:CDM)'
$ws.Range("D3").Value = 0.058
$ws.Range("E3").Value = 0.11
$ws.Range("K3").Value = 23.9
$ws.Range("L3").Value = 0.1293290043290043
$ws.Range("M3").Value = 21
$ws.Range("N3").Value = 0.03745986443096682
$ws.Range("O3").Value = 0.8786610878661089
$ws.Range("P3").Value = 21
$ws.Range("Q3").Value = 0.03745986443096682
$ws.Range("R3").Value = 0.8786610878661089
$ws.Range("U3").Value = 155.6
$ws.Range("V3").Value = 0.2775597574027827
$ws.Range("W3").Value = 0.04317196531791907
$ws.Range("X3").Value = 0.07018276659928611
$ws.Range("Y3").Value = -0.02701080128136704
$ws.Range("Z3").Value = 0.1627906976744186
$ws.Range("AB3").Value = 0.04807136925437845
$ws.Range("AC3").Value = -0.04807136925437845
$ws.Range("AD3").Value = 632.9
$ws.Range("AF3").Value = 632.9
$ws.Range("AG3").Value = 477.3
$ws.Range("AH3").Value = 0.5302890657729368
$ws.Range("AI3").Value = 0.5316253674926501
$ws.Range("AJ3").Value = 0.459870893149629
$ws.Range("AK3").Value = 0.461203981060972

# --- Row 4 ---
$ws.Range("B4").Value = 'Banque Marocaine pour le Commerce et l''Industrie (CBSE:BCI)'
$ws.Range("D4").Value = -0.00911
$ws.Range("E4").Value = -0.108
$ws.Range("K4").Value = 25.4
$ws.Range("L4").Value = 0.1080851063829787
$ws.Range("U4").Value = 111.2
$ws.Range("V4").Value = 0.1169172537062349
$ws.Range("W4").Value = 0.0337631264123355
$ws.Range("X4").Value = 0.07513194827981708
$ws.Range("Y4").Value = -0.04136882186748157
$ws.Range("Z4").Value = 0.1503865253673271
$ws.Range("AB4").Value = 0.04843241772791118
$ws.Range("AC4").Value = -0.04843241772791118
$ws.Range("AD4").Value = 1273.1
$ws.Range("AF4").Value = 1273.1
$ws.Range("AG4").Value = 1161.9
$ws.Range("AH4").Value = 0.5723855768366154
$ws.Range("AI4").Value = 0.6068160152526215
$ws.Range("AJ4").Value = 0.5498816848083293
$ws.Range("AK4").Value = 0.5848097443124622

# --- Row 5 ---
$ws.Range("B5").Value = 'Bank of Africa (CBSE:BOA)'
$ws.Range("D5").Value = 0.00454
$ws.Range("E5").Value = -0.122
$ws.Range("K5").Value = 114.6
$ws.Range("L5").Value = 0.1001223134719553
$ws.Range("U5").Value = 1379.2
$ws.Range("V5").Value = 0.3942486350513106
$ws.Range("W5").Value = 0.05458182510954467
$ws.Range("X5").Value = 0.1051358932374763
$ws.Range("Y5").Value = -0.05055406812793164
$ws.Range("Z5").Value = 0.141336560308209
$ws.Range("AB5").Value = 0.04972356122208303
$ws.Range("AC5").Value = -0.04972356122208303
$ws.Range("AD5").Value = 9127.6
$ws.Range("AF5").Value = 9127.6
$ws.Range("AG5").Value = 7748.400000000001
$ws.Range("AH5").Value = 0.72292668245432
$ws.Range("AI5").Value = 0.750692908075566
$ws.Range("AJ5").Value = 0.6889487583024354
$ws.Range("AK5").Value = 0.7187955137898087

# --- Row 6 ---
$ws.Range("B6").Value = 'Attijariwafa bank S.A (CBSE:ATW)'
$ws.Range("D6").Value = 0.0267
$ws.Range("E6").Value = -0.0532
$ws.Range("K6").Value = 366.9
$ws.Range("L6").Value = 0.1766490129995185
$ws.Range("U6").Value = 2709.3
$ws.Range("V6").Value = 0.2899849083261086
$ws.Range("W6").Value = 0.07678623749529111
$ws.Range("X6").Value = 0.05416648736240148
$ws.Range("Y6").Value = 0.02261975013288962
$ws.Range("Z6").Value = 0.2108115789046323
$ws.Range("AB6").Value = 0.05000373786774416
$ws.Range("AC6").Value = -0.05000373786774416
$ws.Range("AD6").Value = 4211
$ws.Range("AF6").Value = 4211
$ws.Range("AG6").Value = 1501.7
$ws.Range("AH6").Value = 0.3106854853584577
$ws.Range("AI6").Value = 0.4119666982987174
$ws.Range("AJ6").Value = 0.1384744481124246
$ws.Range("AK6").Value = 0.1998961716628507

# --- Row 7 ---
$ws.Range("B7").Value = 'Crédit Immobilier et Hôtelier, Société Anonyme (CBSE:CIH)'
$ws.Range("D7").Value = 0.07339999999999999
$ws.Range("E7").Value = -0.0455
$ws.Range("K7").Value = 33.1
$ws.Range("L7").Value = 0.1422432316287065
$ws.Range("M7").Value = 40.752
$ws.Range("N7").Value = 0.05026147015293538
$ws.Range("O7").Value = 1.231178247734139
$ws.Range("P7").Value = 40.752
$ws.Range("Q7").Value = 0.05026147015293538
$ws.Range("R7").Value = 1.231178247734139
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 80.7
$ws.Range("V7").Value = 0.09953132708436113
$ws.Range("W7").Value = 0.07147484344633988
$ws.Range("X7").Value = 0.1207814689668619
$ws.Range("Y7").Value = -0.04930662552052202
$ws.Range("Z7").Value = 0.09566683111330375
$ws.Range("AB7").Value = 0.05009214303598859
$ws.Range("AC7").Value = -0.05009214303598859
$ws.Range("AD7").Value = 2652.7
$ws.Range("AF7").Value = 2652.7
$ws.Range("AG7").Value = 2572
$ws.Range("AH7").Value = 0.7659015446802367
$ws.Range("AI7").Value = 0.8186587661636269
$ws.Range("AJ7").Value = 0.7603168972448858
$ws.Range("AK7").Value = 0.8140270920369668

# --- Row 8 ---
$ws.Range("B8").Value = 'Banque Centrale Populaire (CBSE:BCP)'
$ws.Range("D8").Value = 0.0383
$ws.Range("E8").Value = 0.0385
$ws.Range("K8").Value = 244.7
$ws.Range("L8").Value = 0.1597258485639687
$ws.Range("M8").Value = 183.2325
$ws.Range("N8").Value = 0.02942737609610381
$ws.Range("O8").Value = 0.7488046587658357
$ws.Range("P8").Value = 183.2325
$ws.Range("Q8").Value = 0.02942737609610381
$ws.Range("R8").Value = 0.7488046587658357
$ws.Range("T8").Value = 0
$ws.Range("U8").Value = 1766.9
$ws.Range("V8").Value = 0.2837664214820287
$ws.Range("W8").Value = 0.07032215420869616
$ws.Range("X8").Value = 0.06539086715093009
$ws.Range("Y8").Value = 0.004931287057766062
$ws.Range("Z8").Value = 0.1915239404925616
$ws.Range("AB8").Value = 0.05355211988957206
$ws.Range("AC8").Value = -0.05355211988957206
$ws.Range("AD8").Value = 5766.1
$ws.Range("AF8").Value = 5766.1
$ws.Range("AG8").Value = 3999.2
$ws.Range("AH8").Value = 0.480800820499137
$ws.Range("AI8").Value = 0.5702855334342145
$ws.Range("AJ8").Value = 0.3910892057345147
$ws.Range("AK8").Value = 0.4792905081495686
